$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1949.9155  # H15: 1973.9636 -> 1949.9155
$ws.Cells.Item(15, 9).Value = 1949.9155  # I15: 1973.9636 -> 1949.9155
$ws.Cells.Item(15, 11).Value = 5849.7465  # K15: 5921.8908 -> 5849.7465
$ws.Cells.Item(15, 13).Value = -5680.7465  # M15: -5752.8908 -> -5680.7465
$ws.Cells.Item(17, 8).Value = 2505.617  # H17: 2591.6445 -> 2505.617
$ws.Cells.Item(17, 10).Value = 2603.8667  # J17: 2698.465 -> 2603.8667
$ws.Cells.Item(17, 12).Value = 7811.6001  # L17: 8095.395 -> 7811.6001
$ws.Cells.Item(17, 14).Value = -8147.6001  # N17: -8431.395 -> -8147.6001
$ws.Cells.Item(33, 8).Value = 143.31818  # H33: 160.75 -> 143.31818
$ws.Cells.Item(33, 9).Value = 112.57895  # I33: 134.3125 -> 112.57895
$ws.Cells.Item(33, 10).Value = 338  # J33: 266.5 -> 338
$ws.Cells.Item(33, 11).Value = 112.57895  # K33: 134.3125 -> 112.57895
$ws.Cells.Item(33, 12).Value = 338  # L33: 266.5 -> 338
$ws.Cells.Item(33, 13).Value = 116.42105  # M33: 94.6875 -> 116.42105
$ws.Cells.Item(33, 14).Value = -796  # N33: -724.5 -> -796
$ws.Cells.Item(106, 8).Value = 3228.8  # H106: 3480.4 -> 3228.8
$ws.Cells.Item(106, 9).Value = 3471.1428  # I106: 3828.2856 -> 3471.1428
$ws.Cells.Item(106, 10).Value = 2663.3333  # J106: 2668.6667 -> 2663.3333
$ws.Cells.Item(106, 11).Value = 3471.1428  # K106: 3828.2856 -> 3471.1428
$ws.Cells.Item(106, 12).Value = 2663.3333  # L106: 2668.6667 -> 2663.3333
$ws.Cells.Item(106, 13).Value = -2840.1428  # M106: -3197.2856 -> -2840.1428
$ws.Cells.Item(106, 14).Value = -3925.3333  # N106: -3930.6667 -> -3925.3333
$ws.Cells.Item(121, 8).Value = 3025  # H121: 4250 -> 3025
$ws.Cells.Item(121, 9).Value = 1333.3334  # I121: 400 -> 1333.3334
$ws.Cells.Item(121, 11).Value = 4000.0002  # K121: 1200 -> 4000.0002
$ws.Cells.Item(121, 13).Value = -2253.0002  # M121: 547 -> -2253.0002
$ws.Cells.Item(129, 8).Value = 1055.1014  # H129: 1046.9296 -> 1055.1014
$ws.Cells.Item(129, 9).Value = 440.33334  # I129: 452.75 -> 440.33334
$ws.Cells.Item(129, 10).Value = 1083.0454  # J129: 1082.403 -> 1083.0454
$ws.Cells.Item(129, 11).Value = 1321.00002  # K129: 1358.25 -> 1321.00002
$ws.Cells.Item(129, 12).Value = 3249.1362  # L129: 3247.209 -> 3249.1362
$ws.Cells.Item(129, 13).Value = 3678.99998  # M129: 3641.75 -> 3678.99998
$ws.Cells.Item(129, 14).Value = -13249.1362  # N129: -13247.209 -> -13249.1362
$ws.Cells.Item(141, 8).Value = 3737.647  # H141: 4704 -> 3737.647
$ws.Cells.Item(141, 9).Value = 2316.25  # I141: 2822.5 -> 2316.25
$ws.Cells.Item(141, 10).Value = 5001.1113  # J141: 7526.25 -> 5001.1113
$ws.Cells.Item(141, 11).Value = 6948.75  # K141: 8467.5 -> 6948.75
$ws.Cells.Item(141, 12).Value = 15003.3339  # L141: 22578.75 -> 15003.3339
$ws.Cells.Item(141, 13).Value = -1768.75  # M141: -3287.5 -> -1768.75
$ws.Cells.Item(141, 14).Value = -25363.3339  # N141: -32938.75 -> -25363.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2589.75  # H2: 4162.2 -> 2589.75
$ws.Cells.Item(2, 9).Value = 2717.7  # I2: 4702.75 -> 2717.7
$ws.Cells.Item(2, 10).Value = 1950  # J2: 2000 -> 1950
$ws.Cells.Item(2, 11).Value = 2717.7  # K2: 4702.75 -> 2717.7
$ws.Cells.Item(2, 12).Value = 1950  # L2: 2000 -> 1950
$ws.Cells.Item(2, 13).Value = -2604.7  # M2: -4589.75 -> -2604.7
$ws.Cells.Item(2, 14).Value = -2176  # N2: -2226 -> -2176
$ws.Cells.Item(61, 8).Value = 10707.462  # H61: 6923 -> 10707.462
$ws.Cells.Item(61, 9).Value = 5443.5625  # I61: 6079.643 -> 5443.5625
$ws.Cells.Item(61, 10).Value = 19129.7  # J61: 7906.9165 -> 19129.7
$ws.Cells.Item(61, 11).Value = 5443.5625  # K61: 6079.643 -> 5443.5625
$ws.Cells.Item(61, 12).Value = 19129.7  # L61: 7906.9165 -> 19129.7
$ws.Cells.Item(61, 13).Value = -5231.5625  # M61: -5867.643 -> -5231.5625
$ws.Cells.Item(61, 14).Value = -19553.7  # N61: -8330.916499999999 -> -19553.7
$ws.Cells.Item(74, 8).Value = 4454.6665  # H74: 4153.7617 -> 4454.6665
$ws.Cells.Item(74, 9).Value = 1918.909  # I74: 1779.7297 -> 1918.909
$ws.Cells.Item(74, 10).Value = 18401.334  # J74: 21721.6 -> 18401.334
$ws.Cells.Item(74, 11).Value = 1918.909  # K74: 1779.7297 -> 1918.909
$ws.Cells.Item(74, 12).Value = 18401.334  # L74: 21721.6 -> 18401.334
$ws.Cells.Item(74, 13).Value = -1044.909  # M74: -905.7297000000001 -> -1044.909
$ws.Cells.Item(74, 14).Value = -20149.334  # N74: -23469.6 -> -20149.334
$ws.Cells.Item(77, 8).Value = 4454.6665  # H77: 4153.7617 -> 4454.6665
$ws.Cells.Item(77, 9).Value = 1918.909  # I77: 1779.7297 -> 1918.909
$ws.Cells.Item(77, 10).Value = 18401.334  # J77: 21721.6 -> 18401.334
$ws.Cells.Item(77, 11).Value = 9594.545  # K77: 8898.648500000001 -> 9594.545
$ws.Cells.Item(77, 12).Value = 92006.67  # L77: 108608 -> 92006.67
$ws.Cells.Item(77, 13).Value = -5226.545  # M77: -4530.648500000001 -> -5226.545
$ws.Cells.Item(77, 14).Value = -100742.67  # N77: -117344 -> -100742.67
$ws.Cells.Item(110, 8).Value = 2051.375  # H110: 2058.7144 -> 2051.375
$ws.Cells.Item(110, 10).Value = 2200  # J110: 2400 -> 2200
$ws.Cells.Item(110, 12).Value = 2200  # L110: 2400 -> 2200
$ws.Cells.Item(110, 14).Value = -6290  # N110: -6490 -> -6290
$ws.Cells.Item(111, 8).Value = 79450  # H111: 78900 -> 79450
$ws.Cells.Item(111, 10).Value = 79450  # J111: 78900 -> 79450
$ws.Cells.Item(111, 12).Value = 79450  # L111: 78900 -> 79450
$ws.Cells.Item(111, 14).Value = -87630  # N111: -87080 -> -87630
$ws.Cells.Item(116, 8).Value = 2589.75  # H116: 4162.2 -> 2589.75
$ws.Cells.Item(116, 9).Value = 2717.7  # I116: 4702.75 -> 2717.7
$ws.Cells.Item(116, 10).Value = 1950  # J116: 2000 -> 1950
$ws.Cells.Item(116, 11).Value = 2717.7  # K116: 4702.75 -> 2717.7
$ws.Cells.Item(116, 12).Value = 1950  # L116: 2000 -> 1950
$ws.Cells.Item(116, 13).Value = -423.6999999999998  # M116: -2408.75 -> -423.6999999999998
$ws.Cells.Item(116, 14).Value = -6538  # N116: -6588 -> -6538
$ws.Cells.Item(122, 8).Value = 2005.8334  # H122: 4466083 -> 2005.8334
$ws.Cells.Item(122, 9).Value = 1800.8  # I122: 1683.3334 -> 1800.8
$ws.Cells.Item(122, 10).Value = 2347.5557  # J122: 17859282 -> 2347.5557
$ws.Cells.Item(122, 11).Value = 5402.4  # K122: 5050.0002 -> 5402.4
$ws.Cells.Item(122, 12).Value = 7042.6671  # L122: 53577846 -> 7042.6671
$ws.Cells.Item(122, 13).Value = -2952.4  # M122: -2600.0002 -> -2952.4
$ws.Cells.Item(122, 14).Value = -11942.6671  # N122: -53582746 -> -11942.6671
$ws.Cells.Item(132, 8).Value = 2157.9778  # H132: 2386.9473 -> 2157.9778
$ws.Cells.Item(132, 9).Value = 1574.2  # I132: 1511.8636 -> 1574.2
$ws.Cells.Item(132, 10).Value = 4201.2  # J132: 3590.1875 -> 4201.2
$ws.Cells.Item(132, 11).Value = 4722.6  # K132: 4535.5908 -> 4722.6
$ws.Cells.Item(132, 12).Value = 12603.6  # L132: 10770.5625 -> 12603.6
$ws.Cells.Item(132, 13).Value = -2192.6  # M132: -2005.5908 -> -2192.6
$ws.Cells.Item(132, 14).Value = -17663.6  # N132: -15830.5625 -> -17663.6
$ws.Cells.Item(136, 8).Value = 10707.462  # H136: 6923 -> 10707.462
$ws.Cells.Item(136, 9).Value = 5443.5625  # I136: 6079.643 -> 5443.5625
$ws.Cells.Item(136, 10).Value = 19129.7  # J136: 7906.9165 -> 19129.7
$ws.Cells.Item(136, 11).Value = 16330.6875  # K136: 18238.929 -> 16330.6875
$ws.Cells.Item(136, 12).Value = 57389.10000000001  # L136: 23720.7495 -> 57389.10000000001
$ws.Cells.Item(136, 13).Value = -13780.6875  # M136: -15688.929 -> -13780.6875
$ws.Cells.Item(136, 14).Value = -62489.10000000001  # N136: -28820.7495 -> -62489.10000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2589.75  # H3: 4162.2 -> 2589.75
$ws.Cells.Item(3, 9).Value = 2717.7  # I3: 4702.75 -> 2717.7
$ws.Cells.Item(3, 10).Value = 1950  # J3: 2000 -> 1950
$ws.Cells.Item(3, 11).Value = 2717.7  # K3: 4702.75 -> 2717.7
$ws.Cells.Item(3, 12).Value = 1950  # L3: 2000 -> 1950
$ws.Cells.Item(3, 13).Value = -2603.7  # M3: -4588.75 -> -2603.7
$ws.Cells.Item(3, 14).Value = -2178  # N3: -2228 -> -2178
$ws.Cells.Item(122, 8).Value = 0  # H122: 44444 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 44444 -> 0
$ws.Cells.Item(122, 12).Value = 0  # L122: 44444 -> 0
$ws.Cells.Item(122, 14).ClearContents()  # N122 was -54244
$ws.Cells.Item(134, 8).Value = 42675.56  # H134: 3219.4783 -> 42675.56
$ws.Cells.Item(134, 9).Value = 2755.4348  # I134: 3332.2778 -> 2755.4348
$ws.Cells.Item(134, 10).Value = 501757  # J134: 2813.4 -> 501757
$ws.Cells.Item(134, 11).Value = 8266.304400000001  # K134: 9996.8334 -> 8266.304400000001
$ws.Cells.Item(134, 12).Value = 1505271  # L134: 8440.200000000001 -> 1505271
$ws.Cells.Item(134, 13).Value = -5731.304400000001  # M134: -7461.8334 -> -5731.304400000001
$ws.Cells.Item(134, 14).Value = -1510341  # N134: -13510.2 -> -1510341

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 139.23529  # H7: 155.78261 -> 139.23529
$ws.Cells.Item(7, 9).Value = 104.78571  # I7: 99.933334 -> 104.78571
$ws.Cells.Item(7, 10).Value = 300  # J7: 260.5 -> 300
$ws.Cells.Item(7, 11).Value = 104.78571  # K7: 99.933334 -> 104.78571
$ws.Cells.Item(7, 12).Value = 300  # L7: 260.5 -> 300
$ws.Cells.Item(7, 13).Value = 8.214290000000005  # M7: 13.066666 -> 8.214290000000005
$ws.Cells.Item(7, 14).Value = -526  # N7: -486.5 -> -526
$ws.Cells.Item(16, 8).Value = 1068.1111  # H16: 1200 -> 1068.1111
$ws.Cells.Item(16, 9).Value = 816.6667  # I16: 1000 -> 816.6667
$ws.Cells.Item(16, 10).Value = 1571  # J16: 2000 -> 1571
$ws.Cells.Item(16, 11).Value = 816.6667  # K16: 1000 -> 816.6667
$ws.Cells.Item(16, 12).Value = 1571  # L16: 2000 -> 1571
$ws.Cells.Item(16, 13).Value = -529.6667  # M16: -713 -> -529.6667
$ws.Cells.Item(16, 14).Value = -2145  # N16: -2574 -> -2145
$ws.Cells.Item(22, 8).Value = 281.875  # H22: 259.44446 -> 281.875
$ws.Cells.Item(22, 10).Value = 298  # J22: 243.5 -> 298
$ws.Cells.Item(22, 12).Value = 298  # L22: 243.5 -> 298
$ws.Cells.Item(22, 14).Value = -998  # N22: -943.5 -> -998
$ws.Cells.Item(31, 8).Value = 2664.7441  # H31: 2585.3333 -> 2664.7441
$ws.Cells.Item(31, 9).Value = 1919.4667  # I31: 1785.1765 -> 1919.4667
$ws.Cells.Item(31, 10).Value = 3064  # J31: 3071.1428 -> 3064
$ws.Cells.Item(31, 11).Value = 1919.4667  # K31: 1785.1765 -> 1919.4667
$ws.Cells.Item(31, 12).Value = 3064  # L31: 3071.1428 -> 3064
$ws.Cells.Item(31, 13).Value = -1624.4667  # M31: -1490.1765 -> -1624.4667
$ws.Cells.Item(31, 14).Value = -3654  # N31: -3661.1428 -> -3654
$ws.Cells.Item(34, 8).Value = 2664.7441  # H34: 2585.3333 -> 2664.7441
$ws.Cells.Item(34, 9).Value = 1919.4667  # I34: 1785.1765 -> 1919.4667
$ws.Cells.Item(34, 10).Value = 3064  # J34: 3071.1428 -> 3064
$ws.Cells.Item(34, 11).Value = 1919.4667  # K34: 1785.1765 -> 1919.4667
$ws.Cells.Item(34, 12).Value = 3064  # L34: 3071.1428 -> 3064
$ws.Cells.Item(34, 13).Value = -1717.4667  # M34: -1583.1765 -> -1717.4667
$ws.Cells.Item(34, 14).Value = -3468  # N34: -3475.1428 -> -3468
$ws.Cells.Item(58, 8).Value = 3369863  # H58: 3369795.8 -> 3369863
$ws.Cells.Item(58, 9).Value = 6495084  # I58: 6062025 -> 6495084
$ws.Cells.Item(58, 10).Value = 4240.3076  # J58: 4509.1665 -> 4240.3076
$ws.Cells.Item(58, 11).Value = 6495084  # K58: 6062025 -> 6495084
$ws.Cells.Item(58, 12).Value = 4240.3076  # L58: 4509.1665 -> 4240.3076
$ws.Cells.Item(58, 13).Value = -6494881  # M58: -6061822 -> -6494881
$ws.Cells.Item(58, 14).Value = -4646.3076  # N58: -4915.1665 -> -4646.3076
$ws.Cells.Item(113, 8).Value = 1068.1111  # H113: 1200 -> 1068.1111
$ws.Cells.Item(113, 9).Value = 816.6667  # I113: 1000 -> 816.6667
$ws.Cells.Item(113, 10).Value = 1571  # J113: 2000 -> 1571
$ws.Cells.Item(113, 11).Value = 816.6667  # K113: 1000 -> 816.6667
$ws.Cells.Item(113, 12).Value = 1571  # L113: 2000 -> 1571
$ws.Cells.Item(113, 13).Value = 1353.3333  # M113: 1170 -> 1353.3333
$ws.Cells.Item(113, 14).Value = -5911  # N113: -6340 -> -5911
$ws.Cells.Item(132, 8).Value = 2408.868  # H132: 2456.1482 -> 2408.868
$ws.Cells.Item(132, 9).Value = 2197.1904  # I132: 2160.1025 -> 2197.1904
$ws.Cells.Item(132, 10).Value = 3217.0908  # J132: 3225.8667 -> 3217.0908
$ws.Cells.Item(132, 11).Value = 6591.5712  # K132: 6480.3075 -> 6591.5712
$ws.Cells.Item(132, 12).Value = 9651.2724  # L132: 9677.6001 -> 9651.2724
$ws.Cells.Item(132, 13).Value = -4061.5712  # M132: -3950.3075 -> -4061.5712
$ws.Cells.Item(132, 14).Value = -14711.2724  # N132: -14737.6001 -> -14711.2724
$ws.Cells.Item(134, 8).Value = 2473.7742  # H134: 2861.9614 -> 2473.7742
$ws.Cells.Item(134, 9).Value = 2416.762  # I134: 2683.8948 -> 2416.762
$ws.Cells.Item(134, 10).Value = 2593.5  # J134: 3345.2856 -> 2593.5
$ws.Cells.Item(134, 11).Value = 7250.286  # K134: 8051.6844 -> 7250.286
$ws.Cells.Item(134, 12).Value = 7780.5  # L134: 10035.8568 -> 7780.5
$ws.Cells.Item(134, 13).Value = -4715.286  # M134: -5516.6844 -> -4715.286
$ws.Cells.Item(134, 14).Value = -12850.5  # N134: -15105.8568 -> -12850.5
$ws.Cells.Item(136, 8).Value = 3369863  # H136: 3369795.8 -> 3369863
$ws.Cells.Item(136, 9).Value = 6495084  # I136: 6062025 -> 6495084
$ws.Cells.Item(136, 10).Value = 4240.3076  # J136: 4509.1665 -> 4240.3076
$ws.Cells.Item(136, 11).Value = 19485252  # K136: 18186075 -> 19485252
$ws.Cells.Item(136, 12).Value = 12720.9228  # L136: 13527.4995 -> 12720.9228
$ws.Cells.Item(136, 13).Value = -19482702  # M136: -18183525 -> -19482702
$ws.Cells.Item(136, 14).Value = -17820.9228  # N136: -18627.4995 -> -17820.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 403.25  # H98: 430.6 -> 403.25
$ws.Cells.Item(98, 9).Value = 386.09525  # I98: 409.4737 -> 386.09525
$ws.Cells.Item(98, 10).Value = 454.7143  # J98: 497.5 -> 454.7143
$ws.Cells.Item(98, 11).Value = 1158.28575  # K98: 1228.4211 -> 1158.28575
$ws.Cells.Item(98, 12).Value = 1364.1429  # L98: 1492.5 -> 1364.1429
$ws.Cells.Item(98, 13).Value = 339.71425  # M98: 269.5789 -> 339.71425
$ws.Cells.Item(98, 14).Value = -4360.1429  # N98: -4488.5 -> -4360.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4265.1904  # H122: 6790.364 -> 4265.1904
$ws.Cells.Item(122, 9).Value = 5427.923  # I122: 10214.333 -> 5427.923
$ws.Cells.Item(122, 10).Value = 2375.75  # J122: 2681.6 -> 2375.75
$ws.Cells.Item(122, 11).Value = 16283.769  # K122: 30642.999 -> 16283.769
$ws.Cells.Item(122, 12).Value = 7127.25  # L122: 8044.799999999999 -> 7127.25
$ws.Cells.Item(122, 13).Value = -13833.769  # M122: -28192.999 -> -13833.769
$ws.Cells.Item(122, 14).Value = -12027.25  # N122: -12944.8 -> -12027.25
$ws.Cells.Item(132, 8).Value = 10810.357  # H132: 3057.4 -> 10810.357
$ws.Cells.Item(132, 9).Value = 3609.2222  # I132: 2921.4546 -> 3609.2222
$ws.Cells.Item(132, 10).Value = 23772.4  # J132: 3223.5557 -> 23772.4
$ws.Cells.Item(132, 11).Value = 10827.6666  # K132: 8764.363799999999 -> 10827.6666
$ws.Cells.Item(132, 12).Value = 71317.20000000001  # L132: 9670.667099999999 -> 71317.20000000001
$ws.Cells.Item(132, 13).Value = -8297.6666  # M132: -6234.363799999999 -> -8297.6666
$ws.Cells.Item(132, 14).Value = -76377.20000000001  # N132: -14730.6671 -> -76377.20000000001
$ws.Cells.Item(137, 8).Value = 0  # H137: 49800 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 49800 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 49800 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137 was -60000

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 929562.6  # H61: 1022428.2 -> 929562.6
$ws.Cells.Item(61, 9).Value = 24798.223  # I61: 31467.428 -> 24798.223
$ws.Cells.Item(61, 10).Value = 5001002.5  # J61: 3334670 -> 5001002.5
$ws.Cells.Item(61, 11).Value = 24798.223  # K61: 31467.428 -> 24798.223
$ws.Cells.Item(61, 12).Value = 5001002.5  # L61: 3334670 -> 5001002.5
$ws.Cells.Item(61, 13).Value = -24596.223  # M61: -31265.428 -> -24596.223
$ws.Cells.Item(61, 14).Value = -5001406.5  # N61: -3335074 -> -5001406.5
$ws.Cells.Item(113, 8).Value = 929562.6  # H113: 1022428.2 -> 929562.6
$ws.Cells.Item(113, 9).Value = 24798.223  # I113: 31467.428 -> 24798.223
$ws.Cells.Item(113, 10).Value = 5001002.5  # J113: 3334670 -> 5001002.5
$ws.Cells.Item(113, 11).Value = 24798.223  # K113: 31467.428 -> 24798.223
$ws.Cells.Item(113, 12).Value = 5001002.5  # L113: 3334670 -> 5001002.5
$ws.Cells.Item(113, 13).Value = -22628.223  # M113: -29297.428 -> -22628.223
$ws.Cells.Item(113, 14).Value = -5005342.5  # N113: -3339010 -> -5005342.5
$ws.Cells.Item(132, 8).Value = 4307.519  # H132: 4528.7754 -> 4307.519
$ws.Cells.Item(132, 9).Value = 4256.5386  # I132: 4365.9736 -> 4256.5386
$ws.Cells.Item(132, 10).Value = 4460.4614  # J132: 5091.1816 -> 4460.4614
$ws.Cells.Item(132, 11).Value = 12769.6158  # K132: 13097.9208 -> 12769.6158
$ws.Cells.Item(132, 12).Value = 13381.3842  # L132: 15273.5448 -> 13381.3842
$ws.Cells.Item(132, 13).Value = -10239.6158  # M132: -10567.9208 -> -10239.6158
$ws.Cells.Item(132, 14).Value = -18441.3842  # N132: -20333.5448 -> -18441.3842
$ws.Cells.Item(136, 8).Value = 4649.6875  # H136: 4286.3774 -> 4649.6875
$ws.Cells.Item(136, 9).Value = 2698.9285  # I136: 2445.7188 -> 2698.9285
$ws.Cells.Item(136, 10).Value = 7380.75  # J136: 7091.1904 -> 7380.75
$ws.Cells.Item(136, 11).Value = 8096.7855  # K136: 7337.1564 -> 8096.7855
$ws.Cells.Item(136, 12).Value = 22142.25  # L136: 21273.5712 -> 22142.25
$ws.Cells.Item(136, 13).Value = -5546.7855  # M136: -4787.1564 -> -5546.7855
$ws.Cells.Item(136, 14).Value = -27242.25  # N136: -26373.5712 -> -27242.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 5315.273  # H113: 3520.9707 -> 5315.273
$ws.Cells.Item(113, 9).Value = 6401.9414  # I113: 4770.7827 -> 6401.9414
$ws.Cells.Item(113, 10).Value = 1620.6  # J113: 907.7273 -> 1620.6
$ws.Cells.Item(113, 11).Value = 19205.8242  # K113: 14312.3481 -> 19205.8242
$ws.Cells.Item(113, 12).Value = 4861.799999999999  # L113: 2723.1819 -> 4861.799999999999
$ws.Cells.Item(113, 13).Value = -17035.8242  # M113: -12142.3481 -> -17035.8242
$ws.Cells.Item(113, 14).Value = -9201.799999999999  # N113: -7063.1819 -> -9201.799999999999
$ws.Cells.Item(132, 8).Value = 1538.4681  # H132: 1559.3469 -> 1538.4681
$ws.Cells.Item(132, 9).Value = 857.875  # I132: 834.9091 -> 857.875
$ws.Cells.Item(132, 10).Value = 2990.4  # J132: 3053.5 -> 2990.4
$ws.Cells.Item(132, 11).Value = 2573.625  # K132: 2504.7273 -> 2573.625
$ws.Cells.Item(132, 12).Value = 8971.200000000001  # L132: 9160.5 -> 8971.200000000001
$ws.Cells.Item(132, 13).Value = -43.625  # M132: 25.27269999999999 -> -43.625
$ws.Cells.Item(132, 14).Value = -14031.2  # N132: -14220.5 -> -14031.2
$ws.Cells.Item(136, 8).Value = 8562.424000000001  # H136: 6155 -> 8562.424000000001
$ws.Cells.Item(136, 9).Value = 7032.5  # I136: 2504.1667 -> 7032.5
$ws.Cells.Item(136, 10).Value = 9436.666999999999  # J136: 9284.286 -> 9436.666999999999
$ws.Cells.Item(136, 11).Value = 21097.5  # K136: 7512.500100000001 -> 21097.5
$ws.Cells.Item(136, 12).Value = 28310.001  # L136: 27852.858 -> 28310.001
$ws.Cells.Item(136, 13).Value = -18547.5  # M136: -4962.500100000001 -> -18547.5
$ws.Cells.Item(136, 14).Value = -33410.001  # N136: -32952.858 -> -33410.001
